$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jalon 2")

# Update durations in rows 3 and 4 (formulas in column G recalc automatically)
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 12

# Fill in row 6 (task "Maquette" / git setup row) with new data.
# G6 becomes a literal string value ("non applicable"), replacing its prior formula.
$ws.Range("C6").Value = "git"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Non applicable"
$ws.Range("F6").Value = "non applicable"
$ws.Range("G6").Value = "non applicable"

# Add new rows 10 and 11 entries
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "git"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "SCSS, CSS, appliquer le template.css à toutes les pages"

# Resize columns to (best)fit the new, wider content
$ws.Columns.Item(2).ColumnWidth = 44.3
$ws.Columns.Item(5).ColumnWidth = 12.3
$ws.Columns.Item(6).ColumnWidth = 12.0

# Update the active selection to D11
$ws.Range("D11").Select()
